# Actualización 11 de Mayo - Tarde
# Insert a new rescatable student row (JUAREZ MORO DENISSE) into the
# "Rescatables" sheet, right before the existing "DE LOS SANTOS XOTLANIHUA
# JENNIFER" row, pushing it down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row above current row 4 (the DE LOS SANTOS... row), which
# shifts that row (and everything below it) down to row 5.
$ws.Rows.Item(4).Insert()

# Fill in the new row 4 with the new rescatable student's data.
$ws.Cells.Item(4, 1).Value = 19330051920201
$ws.Cells.Item(4, 2).Value = "JUAREZ"
$ws.Cells.Item(4, 3).Value = "MORO"
$ws.Cells.Item(4, 4).Value = "DENISSE"
$ws.Cells.Item(4, 5).Value = "ECOLOGÍA"
$ws.Cells.Item(4, 6).Value = "4BLCM"
$ws.Cells.Item(4, 7).Value = 2
